$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2667
$ws.Range("K3").Value = 2575
$ws.Range("K4").Value = 537
$ws.Range("K5").Value = 171
$ws.Range("K6").Value = 3200
$ws.Range("K7").Value = 9150

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 134

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 178
$ws.Range("K6").Value = 202
$ws.Range("K7").Value = 605

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 101
$ws.Range("K3").Value = 130
$ws.Range("K6").Value = 103
$ws.Range("K7").Value = 362

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 101
$ws.Range("K6").Value = 93
$ws.Range("K7").Value = 301

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 64
$ws.Range("K7").Value = 214

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 61
$ws.Range("K5").Value = 4
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K5").Value = 17
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 276
$ws.Range("K8").Value = 605
$ws.Range("K11").Value = 193
$ws.Range("K14").Value = 56
$ws.Range("K15").Value = 89
$ws.Range("K18").Value = 61
$ws.Range("K19").Value = 267
$ws.Range("K20").Value = 209
$ws.Range("K21").Value = 25
$ws.Range("K26").Value = 12
$ws.Range("K29").Value = 474
$ws.Range("K31").Value = 106
$ws.Range("K33").Value = 362
$ws.Range("K36").Value = 106
$ws.Range("K37").Value = 301
$ws.Range("K41").Value = 81
$ws.Range("K42").Value = 316
$ws.Range("K44").Value = 88
$ws.Range("K48").Value = 112
$ws.Range("K51").Value = 100
$ws.Range("K53").Value = 134
$ws.Range("K54").Value = 168
$ws.Range("K55").Value = 101
$ws.Range("K59").Value = 16
$ws.Range("K60").Value = 58
$ws.Range("K63").Value = 35
$ws.Range("K65").Value = 214
$ws.Range("K67").Value = 353
$ws.Range("K71").Value = 28
$ws.Range("K76").Value = 136
$ws.Range("K77").Value = 63
$ws.Range("K78").Value = 127
$ws.Range("K79").Value = 234
$ws.Range("K83").Value = 202
$ws.Range("K84").Value = 65
$ws.Range("K85").Value = 438
$ws.Range("K88").Value = 104
$ws.Range("K89").Value = 120
$ws.Range("K90").Value = 84
$ws.Range("K92").Value = 36
$ws.Range("K95").Value = 146
$ws.Range("K97").Value = 79
$ws.Range("K99").Value = 164
$ws.Range("K101").Value = 9150

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 36
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 112
$ws.Range("K3").Value = 110
$ws.Range("K7").Value = 353

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 127
$ws.Range("K3").Value = 158
$ws.Range("K6").Value = 151
$ws.Range("K7").Value = 474

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 84
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 267

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 18
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 28
$ws.Range("K3").Value = 14
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 80
$ws.Range("K3").Value = 101
$ws.Range("K6").Value = 122
$ws.Range("K7").Value = 316

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 36
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 234

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 73
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 209

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 34
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 92
$ws.Range("K3").Value = 86
$ws.Range("K4").Value = 11
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 276

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 193

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K3").Value = 7
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 159
$ws.Range("K3").Value = 151
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 438

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 63
